# tbm.xlsx revision: "revisione con dati tbm corretti, grafici con testi
# leggibili, utilizzo main thrust force anziche auxiliary"
#
# 1) Workbook calculation mode: manual -> automatic (drops calcMode="manual"
#    from <calcPr>).
# 2) Foglio1!I4 (excavationDiameter row 4): 5.04 -> 5.13
# 3) Foglio1!U4 (nominalThrustForce row 4): 13700 -> 17700  (use MAIN thrust
#    force value here instead of the auxiliary one)
# 4) Add helper rows 7:9 with the "opening ratio" style formulas, one per
#    data row 2/3/4:
#       J(n+5) = J(n) + (I(n)-H(n))/2
#       K(n+5) = 2 * J(n+5)
# 5) Update the view: scroll/select so W4 is the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) switch calculation back to automatic -----------------------------
$excel.Calculation = -4105   # xlCalculationAutomatic

# --- 2) & 3) corrected TBM data values ------------------------------------
$ws.Range("I4").Value = 5.13
$ws.Range("U4").Value = 17700

# --- 4) new opening-ratio helper rows (7, 8, 9) ---------------------------
$ws.Range("J7").Formula = "=J2+(I2-H2)/2"
$ws.Range("K7").Formula = "=2*J7"

$ws.Range("J8").Formula = "=J3+(I3-H3)/2"
$ws.Range("K8").Formula = "=2*J8"

$ws.Range("J9").Formula = "=J4+(I4-H4)/2"
$ws.Range("K9").Formula = "=2*J9"

# --- 5) view/selection state ----------------------------------------------
# Scroll the visible pane so column C is the leftmost visible column, then
# make W4 the active/selected cell (matches the saved sheetView state).
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1

$ws.Range("W4").Select()
